$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulated run metrics for data_set gr25_03 (rows 2-11)
# Values taken from the updated OOXML diff for columns C and E:Y

$ws.Range("C2").Value = 0.4879930019378662
$ws.Range("E2").Value = 75.47239478450501
$ws.Range("F2").Value = 0.004279589145457027
$ws.Range("G2").Value = 0.002972863424041272
$ws.Range("H2").Value = 0.002356498816351619
$ws.Range("I2").Value = 0.002196903151732868
$ws.Range("J2").Value = 0.002196903151732868
$ws.Range("K2").Value = 0.002196903151732868
$ws.Range("L2").Value = 0.002095163540351789
$ws.Range("M2").Value = 0.002042883576967383
$ws.Range("N2").Value = 0.001900139670317435
$ws.Range("O2").Value = 0.001900139670317435
$ws.Range("P2").Value = 0.001900139670317435
$ws.Range("Q2").Value = 0.001900139670317435
$ws.Range("R2").Value = 0.001682157296896693
$ws.Range("S2").Value = 0.001682157296896693
$ws.Range("T2").Value = 0.001572425425249419
$ws.Range("U2").Value = 0.001550411410071621
$ws.Range("V2").Value = 0.001550411410071621
$ws.Range("W2").Value = 0.001502005590991083
$ws.Range("X2").Value = 0.001493583462720318
$ws.Range("Y2").Value = 0.001471196779425049

$ws.Range("C3").Value = 0.5249998569488525
$ws.Range("E3").Value = 79.87253807417983
$ws.Range("F3").Value = 0.004288146682722698
$ws.Range("G3").Value = 0.002961314588823508
$ws.Range("H3").Value = 0.00259884869322035
$ws.Range("I3").Value = 0.002423312668613885
$ws.Range("J3").Value = 0.002423312668613885
$ws.Range("K3").Value = 0.002246510042613024
$ws.Range("L3").Value = 0.001913138130489615
$ws.Range("M3").Value = 0.001913138130489615
$ws.Range("N3").Value = 0.001770737406600142
$ws.Range("O3").Value = 0.001770737406600142
$ws.Range("P3").Value = 0.001682547516583332
$ws.Range("Q3").Value = 0.001682547516583332
$ws.Range("R3").Value = 0.0016156032200068
$ws.Range("S3").Value = 0.0016156032200068
$ws.Range("T3").Value = 0.001593486296888197
$ws.Range("U3").Value = 0.00158399732875918
$ws.Range("V3").Value = 0.001579179392829352
$ws.Range("W3").Value = 0.001565426010752047
$ws.Range("X3").Value = 0.001559024844302388
$ws.Range("Y3").Value = 0.00155696955310292

$ws.Range("C4").Value = 0.3839993476867676
$ws.Range("E4").Value = 81.30939133463289
$ws.Range("F4").Value = 0.004229258877639737
$ws.Range("G4").Value = 0.003342676720869958
$ws.Range("H4").Value = 0.002480200667251429
$ws.Range("I4").Value = 0.002480200667251429
$ws.Range("J4").Value = 0.002480200667251429
$ws.Range("K4").Value = 0.002353892344172531
$ws.Range("L4").Value = 0.00232882162663681
$ws.Range("M4").Value = 0.001961494460421718
$ws.Range("N4").Value = 0.001961494460421718
$ws.Range("O4").Value = 0.001953882973847742
$ws.Range("P4").Value = 0.001903203595835829
$ws.Range("Q4").Value = 0.00188225023550554
$ws.Range("R4").Value = 0.001870307968947812
$ws.Range("S4").Value = 0.001818391634143893
$ws.Range("T4").Value = 0.001709236584771083
$ws.Range("U4").Value = 0.001709236584771083
$ws.Range("V4").Value = 0.001666505883086467
$ws.Range("W4").Value = 0.001637579707364662
$ws.Range("X4").Value = 0.00160883526409487
$ws.Range("Y4").Value = 0.001584978388589335

$ws.Range("C5").Value = 0.3659930229187012
$ws.Range("E5").Value = 78.66525120359984
$ws.Range("F5").Value = 0.004089854668626796
$ws.Range("G5").Value = 0.003226298648492162
$ws.Range("H5").Value = 0.002712361071711034
$ws.Range("I5").Value = 0.002284534816118364
$ws.Range("J5").Value = 0.002174935472780311
$ws.Range("K5").Value = 0.002022385284955186
$ws.Range("L5").Value = 0.001768586136797778
$ws.Range("M5").Value = 0.001717163249246531
$ws.Range("N5").Value = 0.001681251534613601
$ws.Range("O5").Value = 0.001681251534613601
$ws.Range("P5").Value = 0.001681251534613601
$ws.Range("Q5").Value = 0.001621727749781967
$ws.Range("R5").Value = 0.001621727749781967
$ws.Range("S5").Value = 0.001621727749781967
$ws.Range("T5").Value = 0.001600265392164189
$ws.Range("U5").Value = 0.001578415138564485
$ws.Range("V5").Value = 0.001566037571066442
$ws.Range("W5").Value = 0.001563037498513421
$ws.Range("X5").Value = 0.001537763204898588
$ws.Range("Y5").Value = 0.001533435695976605

$ws.Range("C6").Value = 0.3900017738342285
$ws.Range("E6").Value = 80.40645527661945
$ws.Range("F6").Value = 0.00395257805847371
$ws.Range("G6").Value = 0.002869890448008141
$ws.Range("H6").Value = 0.002869890448008141
$ws.Range("I6").Value = 0.002453492079144553
$ws.Range("J6").Value = 0.002307262246005146
$ws.Range("K6").Value = 0.002132370505371698
$ws.Range("L6").Value = 0.001924185823046327
$ws.Range("M6").Value = 0.001856003666718436
$ws.Range("N6").Value = 0.001856003666718436
$ws.Range("O6").Value = 0.001856003666718436
$ws.Range("P6").Value = 0.001842515574934043
$ws.Range("Q6").Value = 0.001830131184793877
$ws.Range("R6").Value = 0.001719671307761536
$ws.Range("S6").Value = 0.001706311772232211
$ws.Range("T6").Value = 0.001659002164157682
$ws.Range("U6").Value = 0.001659002164157682
$ws.Range("V6").Value = 0.001645105340877089
$ws.Range("W6").Value = 0.001615934561807122
$ws.Range("X6").Value = 0.001583811641952916
$ws.Range("Y6").Value = 0.001567377295840535

$ws.Range("C7").Value = 0.3389954566955566
$ws.Range("E7").Value = 84.00453912342891
$ws.Range("F7").Value = 0.004288146682722698
$ws.Range("G7").Value = 0.00286273414262498
$ws.Range("H7").Value = 0.002621704162942988
$ws.Range("I7").Value = 0.002272015642620299
$ws.Range("J7").Value = 0.002098305586200775
$ws.Range("K7").Value = 0.002043731608008629
$ws.Range("L7").Value = 0.001750267856295019
$ws.Range("M7").Value = 0.001750267856295019
$ws.Range("N7").Value = 0.001732031042076527
$ws.Range("O7").Value = 0.001732031042076527
$ws.Range("P7").Value = 0.001732031042076527
$ws.Range("Q7").Value = 0.001732031042076527
$ws.Range("R7").Value = 0.001725424843201867
$ws.Range("S7").Value = 0.001725424843201867
$ws.Range("T7").Value = 0.001722310827168101
$ws.Range("U7").Value = 0.001685242720063912
$ws.Range("V7").Value = 0.001672789801775498
$ws.Range("W7").Value = 0.001655289970399147
$ws.Range("X7").Value = 0.001647010135888581
$ws.Range("Y7").Value = 0.00163751538252298

$ws.Range("C8").Value = 0.3300042152404785
$ws.Range("E8").Value = 80.69385872648127
$ws.Range("F8").Value = 0.004003341928958737
$ws.Range("G8").Value = 0.003115575822753173
$ws.Range("H8").Value = 0.002528414367978852
$ws.Range("I8").Value = 0.002347888577006214
$ws.Range("J8").Value = 0.002257938725685206
$ws.Range("K8").Value = 0.002102911413947697
$ws.Range("L8").Value = 0.002102071701641475
$ws.Range("M8").Value = 0.00204320694280386
$ws.Range("N8").Value = 0.001933941226160203
$ws.Range("O8").Value = 0.001916488846147604
$ws.Range("P8").Value = 0.001703299306794925
$ws.Range("Q8").Value = 0.001699779950309431
$ws.Range("R8").Value = 0.001699779950309431
$ws.Range("S8").Value = 0.001639753777054052
$ws.Range("T8").Value = 0.001638451200585557
$ws.Range("U8").Value = 0.001617845303297297
$ws.Range("V8").Value = 0.001600277771013692
$ws.Range("W8").Value = 0.001592937948735871
$ws.Range("X8").Value = 0.00158230785105609
$ws.Range("Y8").Value = 0.00157297970227059

$ws.Range("C9").Value = 0.4199953079223633
$ws.Range("E9").Value = 81.27112977465185
$ws.Range("F9").Value = 0.003852812270059351
$ws.Range("G9").Value = 0.003109075006972059
$ws.Range("H9").Value = 0.002462799916452727
$ws.Range("I9").Value = 0.002306601379005214
$ws.Range("J9").Value = 0.002250743729020211
$ws.Range("K9").Value = 0.002221761473015098
$ws.Range("L9").Value = 0.002103156410040811
$ws.Range("M9").Value = 0.002103156410040811
$ws.Range("N9").Value = 0.002103156410040811
$ws.Range("O9").Value = 0.002103156410040811
$ws.Range("P9").Value = 0.002019959851137105
$ws.Range("Q9").Value = 0.001951134997145937
$ws.Range("R9").Value = 0.001880008650122616
$ws.Range("S9").Value = 0.00175090749297853
$ws.Range("T9").Value = 0.001744102459913734
$ws.Range("U9").Value = 0.001672894978876508
$ws.Range("V9").Value = 0.001668802419370226
$ws.Range("W9").Value = 0.001641286386869017
$ws.Range("X9").Value = 0.001584911750709665
$ws.Range("Y9").Value = 0.001584232549213486

$ws.Range("C10").Value = 0.4850020408630371
$ws.Range("E10").Value = 80.75011687163533
$ws.Range("F10").Value = 0.004192849859966024
$ws.Range("G10").Value = 0.003152893982304248
$ws.Range("H10").Value = 0.002924202128650289
$ws.Range("I10").Value = 0.002544878685107478
$ws.Range("J10").Value = 0.002138551343546514
$ws.Range("K10").Value = 0.002138551343546514
$ws.Range("L10").Value = 0.002138551343546514
$ws.Range("M10").Value = 0.001901194030747629
$ws.Range("N10").Value = 0.001901194030747629
$ws.Range("O10").Value = 0.001863372389907383
$ws.Range("P10").Value = 0.001811611992414763
$ws.Range("Q10").Value = 0.001791381833890729
$ws.Range("R10").Value = 0.001747643140686591
$ws.Range("S10").Value = 0.001665524231496439
$ws.Range("T10").Value = 0.001665524231496439
$ws.Range("U10").Value = 0.001623000525289621
$ws.Range("V10").Value = 0.001590211976360335
$ws.Range("W10").Value = 0.001583359212431477
$ws.Range("X10").Value = 0.001576893198723508
$ws.Range("Y10").Value = 0.001574076352273593

$ws.Range("C11").Value = 0.4519975185394287
$ws.Range("E11").Value = 68.58527101840627
$ws.Range("F11").Value = 0.004219412379640493
$ws.Range("G11").Value = 0.002940037068640565
$ws.Range("H11").Value = 0.002427784149640962
$ws.Range("I11").Value = 0.002282630553810187
$ws.Range("J11").Value = 0.001992220312447785
$ws.Range("K11").Value = 0.001992220312447785
$ws.Range("L11").Value = 0.001992220312447785
$ws.Range("M11").Value = 0.001880759086960497
$ws.Range("N11").Value = 0.001855474186040648
$ws.Range("O11").Value = 0.00165072724440749
$ws.Range("P11").Value = 0.00165072724440749
$ws.Range("Q11").Value = 0.001641769004731754
$ws.Range("R11").Value = 0.001629436826836657
$ws.Range("S11").Value = 0.00151714131197445
$ws.Range("T11").Value = 0.001458855479752166
$ws.Range("U11").Value = 0.001434413711092887
$ws.Range("V11").Value = 0.001421519844227279
$ws.Range("W11").Value = 0.001408993949470292
$ws.Range("X11").Value = 0.001374429560295695
$ws.Range("Y11").Value = 0.001336944854159966
